$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SIU")

# Rows 4-6: fill in the previously blank "Dia3" scores (columns D, I, N, S, X)
# for each of the five criteria blocks (Eficiencia, Eficacia, Iniciativa,
# Respeto, Puntualidad) with a 10, which ripples through the AVERAGE()
# formulas in F/K/P/U/Z, the Resumen columns AA:AE and the AF Total Mensual.
$ws.Range("D4").Value = 10
$ws.Range("I4").Value = 10
$ws.Range("N4").Value = 10
$ws.Range("S4").Value = 10
$ws.Range("X4").Value = 10

$ws.Range("D5").Value = 10
$ws.Range("I5").Value = 10
$ws.Range("N5").Value = 10
$ws.Range("S5").Value = 10
$ws.Range("X5").Value = 10

$ws.Range("D6").Value = 10
$ws.Range("I6").Value = 10
$ws.Range("N6").Value = 10
$ws.Range("S6").Value = 10
$ws.Range("X6").Value = 10

# Update the saved view state: scroll the sheet back to the top-left
# (clears the stale "topLeftCell" anchor) and leave the selection on X7.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("X7").Select()

$wb.Save()
